$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 72; this shifts the existing rows 72:115 down to 73:116
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new weekly record
$ws.Range("A72").Value = 7
$ws.Range("B72").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C72").Value = "Ñuble"
$ws.Range("D72").Value = 45236
$ws.Range("E72").Value = 16
$ws.Range("F72").Value = 100112022
$ws.Range("G72").Value = "Arveja Verde"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 80
$ws.Range("K72").Value = 22000
$ws.Range("L72").Value = 24000
$ws.Range("M72").Value = 23250
$ws.Range("N72").Value = "$/saco 25 kilos"
$ws.Range("O72").Value = "Región del Maule"
$ws.Range("P72").Value = 930
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"
